$d = $word.ActiveDocument

# Replace the area figure for Alaska (266,420) with the new value (3,937),
# keeping the surrounding spaces intact.
$d.Content.Find.Execute("266,420", $true, $false, $false, $false, $false,
                         $true, 1, $false, "3,937", 2)
